$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (name, total_registros) in final row order (rows 2-11),
# sorted descending by total_registros as in the target workbook.
$data = @(
    @("MAZA RIOFRIO CINTHIA NATELAHI", 129),
    @("VEGA ZAPATA JESUS GABRIEL", 125),
    @("PANTA VARONA CANDY ELIZABETH", 123),
    @("PALACIOS PANTA LUIS MIGUEL", 123),
    @("PANTA NIMA FREDDY ROLAND JUNIOR", 123),
    @("CRISANTO CARMEN ROSITA ABIGAIL", 119),
    @("ELIAS MACHADO JUANA MARGOT", 119),
    @("SALAZAR VEGA MARIA FERNANDA", 119),
    @("HIDALGO MOSCOL YESSICA JAZMIN", 112),
    @("TALLEDO ELIAS ANDREA ALESSANDRA", 107)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
